$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column B (shifting existing B -> D and C -> E)
$ws.Range("B1:C1").EntireColumn.Insert()

# New header values for the newly inserted date columns
# (write C1 first so "Jun_15" is registered in the shared-string table before "Jun_17")
$ws.Range("C1").Value = "Jun_15"
$ws.Range("B1").Value = "Jun_17"

# Fill the new B and C columns (rows 2-27) with "UN" to match the rest of the table
$ws.Range("B2:C27").Value = "UN"

# Re-apply the custom column width (8.0 characters) to columns C, D and E - the
# insert operation shifted the original column C (and its custom width) to E,
# so give the newly inserted C and D columns a matching custom width as well.
$ws.Range("C1:E1").EntireColumn.ColumnWidth = 7.166666666666667
